{"js": "// The two logo pictures live in the primary header as floating (anchored)\n// shapes, not inline pictures, so they are reached via\n// section.getHeader(\"Primary\").shapes rather than body.inlinePictures.\nconst sections = context.document.sections;\nsections.load(\"items\");\nawait context.sync();\n\nconst section = sections.items[0];\nconst header = section.getHeader(\"Primary\");\nconst shapes = header.shapes;\nshapes.load(\"items/name\");\nawait context.sync();\n\nfor (const shape of shapes.items) {\n  if (shape.name === \"image2.png\") {\n    shape.name = \"image1.png\";\n  } else if (shape.name === \"image1.jpg\") {\n    shape.name = \"image2.jpg\";\n  }\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The two logo pictures live in the primary header as floating (anchored)\n# shapes, not inline pictures, so they are reached via\n# Sections(1).Headers(wdHeaderFooterPrimary).Shapes rather than InlineShapes.\n$sec = $d.Sections.Item(1)\n$hdr = $sec.Headers.Item(1)  # wdHeaderFooterPrimary\n$shapes = $hdr.Shapes\n\nfor ($i = 1; $i -le $shapes.Count; $i++) {\n    $s = $shapes.Item($i)\n    if ($s.Name -eq \"image2.png\") {\n        $s.Name = \"image1.png\"\n    } elseif ($s.Name -eq \"image1.jpg\") {\n        $s.Name = \"image2.jpg\"\n    }\n}\n"}
